# Update column F (dSF) values for specific rows to reflect the
# repulled / recalculated data (mean calculation) as per commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    7  = -4
    10 = 0
    15 = -1
    21 = -4
    26 = -3
    30 = 1
    35 = -3
    46 = 3
    47 = 4
    50 = 4
    53 = 1
    54 = 2
    55 = 0
    56 = -2
    66 = -1
    70 = 3
    74 = -1
    75 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
